$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B15").Value = "Gaizka"
$ws.Range("C15").Value = "Extra: Empleados no disponibles del listado de empleados"
$ws.Range("D15").NumberFormat = $ws.Range("D12").NumberFormat
$ws.Range("D15").Value = "5/6/2025"

$ws.Range("D16").Select()
